# Auto-generated Excel COM-interop script to apply impeachment-support edits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5, 7, 11, 24, 25, 27, 41, 42, 43, 67, 70, 83, 89, 111, 119, 129, 130, 135, 145, 155, 162, 167, 170, 183, 191, 196, 209, 214)

foreach ($r in $rows) {
    $ws.Range("A$r").Value = "YES"
    $ws.Range("H$r").Value = 43732
    $ws.Range("H$r").NumberFormat = "yyyy-mm-dd"
    $ws.Range("J$r").Value = 43732
    $ws.Range("J$r").NumberFormat = "yyyy-mm-dd"
    $ws.Range("K$r").Value = 9
    $ws.Range("L$r").Value = 2019
}

$blurbs = @{}
$blurbs[5] = """Beyond an abuse of power, allegations that the President threatened to leverage U.S. taxpayer dollars to extort a foreign government, if true, constitute an unequivocal violation of our federal laws and the U.S. Constitution I swore to defend,"" she tweeted on September 24. ""For the sake of our nat'l security and our democracy, these serious allegations require independent Congressional investigation unobstructed by this Administration. Congress has a responsibility to uphold the rule of law & to take appropriate steps to open an impeachment inquiry."""
$blurbs[7] = """Ever since Trump set foot in the Oval Office, he has been on a mission to rip apart decades of policies that protect civil rights, the environment, public lands and more. He has embarrassed our nation in the eyes of the world with his regular display of ignorance of treaties and historic alliances along with his horrific treatment of families and violation of international laws at our borders,"" she said in a statement on September 24. ""It is time for a formal impeachment inquiry against this lawless and unstable President."""
$blurbs[11] = """Members of Congress took an oath before God to uphold the Constitution, not to the President,"" he said in a news release on September 24. ""It is imperative to our national security, the preservation of the rule of law, and our patriotism that we begin an impeachment inquiry."""
$blurbs[24] = """We need to get to the bottom of these serious allegations through comprehensive impeachment proceedings and document production,"" he tweeted on September 24. ""Whether such hearings result in actual impeachment of the president and/or senior cabinet officials remains to be seen."""
$blurbs[25] = """Together with other claimed breaches and the President's continued obstruction of Congressional oversight, they fully justify Congress reviewing potential impeachment,"" he said in a news release on September 24."
$blurbs[27] = """Trump has violated his oath of office. He operates for himself, in his personal interest -- and not in the interest of the American people,"" she said in a statement on September 24. ""The impeachment inquiry should be put on the fast track and all of the facts made plain."""
$blurbs[41] = """I support an official impeachment investigation,"" he said in a news release on September 24. ""It is paramount the six committees that have jurisdiction continue to investigate whether any high crimes and misdemeanors have been committed by President Trump."""
$blurbs[42] = """We are at a critical moment for our nation, involving a 'credible' complaint concerning the actions of a sitting president happening in real-time,"" he said in a news release on September 24. ""That is why I now believe it is time to elevate this process to a formal inquiry on the President and his potential misconduct, and to follow the facts wherever they may lead -- which includes the most serious action the House can take under the Constitution: impeachment."""
$blurbs[43] = """When the President of the United States is putting his personal and political interests ahead of the American people and abusing the Constitution, we have a duty to investigate and follow the facts where they lead,"" Cox tweeted on September 24. ""An impeachment inquiry will allow us to break through the President's stonewalling and obstruction and get the answers the American people deserve."""
$blurbs[67] = """I fully support the Speaker's announcement of a new impeachment inquiry,"" she tweeted on September 24. ""The circumstances surrounding the President's phone call with a foreign leader requires it and we are obligated to go wherever the facts lead."""
$blurbs[70] = """I have always supported Congress and our House committees in defending our democracy and Constitution. The impeachment inquiry requested today is needed to continue that work,"" she said in a news release on September 24. ""It's imperative we stand up for our country and hold those accountable who hurt our democracy."""
$blurbs[83] = """The serious allegations that have surfaced about communication between President Trump and the Ukrainian President put our national security at risk and merit an immediate inquiry,"" he tweeted on September 24. ""We must let the facts guide our work."""
$blurbs[89] = """The actions of this President lead me to no other conclusion than to support formal impeachment proceedings,"" she said in a news release on September 24."
$blurbs[111] = "Kim said in a statement on September 24, ""By urging the Ukrainian government to take action to influence our Democracy, Trump has violated that power and the trust of the American people. If the facts are corroborated, that violation, and my understanding of its implications, has led me to come to the conclusion that the President has committed an impeachable offense."""
$blurbs[119] = """I have spoken with Rep. Jim Himes, a senior member of the Intelligence Committee, at length who has underscored the gravity of this situation, which is why at today's caucus I will be advocating for impeachment proceedings to begin,"" Larson said in a statement on September 24."
$blurbs[129] = """I have felt the need to proceed judiciously given the high constitutional requirements of impeachment, but President Trump's actions and obstruction require this step of establishing an impeachment inquiry,"" he said in a news release on September 24. ""It is unfortunate that we have ended up at this point."""
$blurbs[130] = """A formal impeachment inquiry should be undertaken in the House of Representatives,"" she said in a news release on September 24. ""The President has now admitted that he asked the President of Ukraine to take actions to help his re-election. That conduct alone violates his obligations under the Constitution. These facts pose a threat to our national security."""
$blurbs[135] = """President Trump's actions, currently undisputed, are an attack on our Democracy, national security and rule of law, and warrant the commencement of formal impeachment proceeding,"" Lynch tweeted on September 24."
$blurbs[145] = """To fulfill our sworn duty, Congress must move forward with an impeachment inquiry. Brazen lies have been a hallmark of this presidency and it is now time for a reckoning,"" he tweeted on September 24."
$blurbs[155] = """I strongly back Speaker of the House Nancy Pelosi's call today for a formal impeachment inquiry into President Donald Trump,"" he tweeted on September 24. ""The American people expect truth and transparency from their government."""
$blurbs[162] = """I support a formal impeachment inquiry into President Donald Trump,"" Panetta tweeted on September 24. ""This inquiry is about uncovering evidence, seeking the truth, and protecting our national security and the integrity of our democracy."""
$blurbs[167] = """I have and continue to support the impeachment investigation,"" Perlmutter tweeted on September 24. ""President Trump needs to be held accountable and he will be."""
$blurbs[170] = """What the President has admitted to already in my estimation is impeachable,"" Phillips told CNN's Erin Burnett on September 24, also calling Pelosi's announcement an ""important step."""
$blurbs[183] = """I am deeply troubled by the recent reports detailing President Trump's pressuring of the Ukrainian government to interfere in our elections. The fact that the President himself confirmed these reports adds to the gravity of the situation,"" he wrote in a Facebook post on September 24. ""This is why, out of respect for our Constitution, concern for our national security, and in the interest of my constituents, I am fully supportive of an official impeachment inquiry."""
$blurbs[191] = """It's bad enough Trump sought help from a foreign power in the last election. It's worse still that he obstructed the investigation into his misconduct. Now he's admitted using his office to coerce another country to interfere in 2020,"" Schiff, the chairman of the House Intelligence Committee, tweeted on September 24. ""I fully support the impeachment inquiry."""
$blurbs[196] = """Every Member of Congress swore an oath to uphold and defend the Constitution. I therefore support Speaker Pelosi's call for a formal impeachment inquiry,"" he said in a news release on September 24. ""It is the only appropriate course of action given the severity of recent events."""
$blurbs[209] = """President Trump may have used the power of his office to pressure a foreign head of state for his own personal political gain. If true, these actions represent an impeachable offense,"" she tweeted on September 24, also issuing a statement calling for ""an impeachment inquiry."""
$blurbs[214] = """The President asking for help on his campaign from a foreign government is against the law with or without a quid pro quo. This is a matter of national security and I fully support Speaker Pelosi's move to an impeachment inquiry,"" Thomson tweeted on September 24."

foreach ($r in $blurbs.Keys) {
    $ws.Range("AL$r").Value = $blurbs[$r]
}

Write-Host "Done applying impeachment edits to $($rows.Count) rows"
